# Updated simulation-derived transition probabilities for the "Saint Peter's_A" matrix sheet
# (commit: added more games, sped up simulate game logic, and drafted optimization logic)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1921824104234528
$ws.Range("C2").Value = 0.5472312703583062
$ws.Range("J2").Value = 0.01954397394136808
$ws.Range("P2").Value = 0.1628664495114006
$ws.Range("S2").Value = 0.0781758957654723

# Row 3
$ws.Range("B3").Value = 0.01754385964912281
$ws.Range("C3").Value = 0.01754385964912281
$ws.Range("J3").Value = 0.03508771929824561
$ws.Range("P3").Value = 0.7368421052631579
$ws.Range("S3").Value = 0.1929824561403509

# Row 4
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.6341463414634146
$ws.Range("S4").Value = 0.2926829268292683

# Row 6
$ws.Range("B6").Value = 0.07407407407407407
$ws.Range("D6").Value = 0.02314814814814815
$ws.Range("F6").Value = 0.0787037037037037
$ws.Range("J6").Value = 0.2638888888888889
$ws.Range("O6").Value = 0.01388888888888889
$ws.Range("Q6").Value = 0.1018518518518518
$ws.Range("R6").Value = 0.07407407407407407

# Row 7
$ws.Range("B7").Value = 0.1282051282051282
$ws.Range("D7").Value = 0.005128205128205128
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1435897435897436
$ws.Range("O7").Value = 0.01025641025641026
$ws.Range("Q7").Value = 0.1487179487179487
$ws.Range("R7").Value = 0.06153846153846154
$ws.Range("S7").Value = 0.4358974358974359

# Row 8
$ws.Range("B8").Value = 0.1220930232558139
$ws.Range("D8").Value = 0.01744186046511628
$ws.Range("F8").Value = 0.04941860465116279
$ws.Range("J8").Value = 0.1540697674418605
$ws.Range("O8").Value = 0.02616279069767442
$ws.Range("Q8").Value = 0.1220930232558139
$ws.Range("R8").Value = 0.0872093023255814
$ws.Range("S8").Value = 0.4215116279069768

# Row 9
$ws.Range("B9").Value = 0.147887323943662
$ws.Range("D9").Value = 0.007042253521126761
$ws.Range("F9").Value = 0.07042253521126761
$ws.Range("J9").Value = 0.1126760563380282
$ws.Range("O9").Value = 0.007042253521126761
$ws.Range("Q9").Value = 0.1267605633802817
$ws.Range("R9").Value = 0.09154929577464789
$ws.Range("S9").Value = 0.4366197183098591

# Row 10
$ws.Range("B10").Value = 0.1148086522462562
$ws.Range("D10").Value = 0.02412645590682196
$ws.Range("F10").Value = 0.08319467554076539
$ws.Range("J10").Value = 0.1480865224625624
$ws.Range("O10").Value = 0.01164725457570715
$ws.Range("Q10").Value = 0.2004991680532446
$ws.Range("R10").Value = 0.05990016638935108
$ws.Range("S10").Value = 0.3577371048252912

# Row 11
$ws.Range("G11").Value = 0.1578947368421053
$ws.Range("J11").Value = 0.1021671826625387
$ws.Range("K11").Value = 0.195046439628483
$ws.Range("L11").Value = 0.5170278637770898
$ws.Range("S11").Value = 0.02786377708978328

# Row 12
$ws.Range("G12").Value = 0.7192982456140351
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("L12").Value = 0.01754385964912281
$ws.Range("S12").Value = 0.04093567251461988

# Row 13
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("S13").Value = 0.02857142857142857

# Row 14
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333

# Row 15
$ws.Range("F15").Value = 0.03571428571428571
$ws.Range("H15").Value = 0.1339285714285714
$ws.Range("I15").Value = 0.04464285714285714
$ws.Range("J15").Value = 0.4330357142857143
$ws.Range("K15").Value = 0.08482142857142858
$ws.Range("M15").Value = 0.01339285714285714
$ws.Range("O15").Value = 0.08928571428571429
$ws.Range("S15").Value = 0.1651785714285714

# Row 16
$ws.Range("F16").Value = 0.005208333333333333
$ws.Range("H16").Value = 0.1041666666666667
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.515625
$ws.Range("K16").Value = 0.1041666666666667
$ws.Range("M16").Value = 0.015625
$ws.Range("O16").Value = 0.06770833333333333
$ws.Range("S16").Value = 0.1041666666666667

# Row 17
$ws.Range("F17").Value = 0.02305475504322766
$ws.Range("H17").Value = 0.2132564841498559
$ws.Range("I17").Value = 0.069164265129683
$ws.Range("J17").Value = 0.3948126801152738
$ws.Range("K17").Value = 0.1296829971181556
$ws.Range("M17").Value = 0.008645533141210375
$ws.Range("N17").Value = 0.005763688760806916
$ws.Range("O17").Value = 0.069164265129683
$ws.Range("S17").Value = 0.08645533141210375

# Row 18
$ws.Range("F18").Value = 0.03496503496503497
$ws.Range("H18").Value = 0.1888111888111888
$ws.Range("I18").Value = 0.06993006993006994
$ws.Range("J18").Value = 0.3566433566433567
$ws.Range("K18").Value = 0.1258741258741259
$ws.Range("M18").Value = 0.02097902097902098
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.1258741258741259

# Row 19
$ws.Range("F19").Value = 0.01702508960573477
$ws.Range("H19").Value = 0.1756272401433692
$ws.Range("I19").Value = 0.07168458781362007
$ws.Range("J19").Value = 0.3655913978494624
$ws.Range("K19").Value = 0.1370967741935484
$ws.Range("M19").Value = 0.02150537634408602
$ws.Range("N19").Value = 0.0008960573476702509
$ws.Range("O19").Value = 0.08870967741935484
$ws.Range("S19").Value = 0.1218637992831541

